$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("G2").Value = 0.4
$ws.Range("E3").Value = 1.26
$ws.Range("B4").Value = 1.43
$ws.Range("C4").Value = 1.41
$ws.Range("E4").Value = 1.24
$ws.Range("F4").Value = 1.07
$ws.Range("C5").Value = 1.38
$ws.Range("D5").Value = 1.32
$ws.Range("F5").Value = 1.04
$ws.Range("G5").Value = 0.7
$ws.Range("D6").Value = 1.55
$ws.Range("G6").Value = 1.05
$ws.Range("B7").Value = 2.6
$ws.Range("F7").Value = 1.49
